$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.875026226043701
$ws.Range("B1").Value = 2.759056091308594
$ws.Range("C1").Value = 2.185009717941284
$ws.Range("D1").Value = 2.045663356781006
$ws.Range("E1").Value = 1.769544720649719
